# Auto-generated: restore swapped/shuffled row data in "Denmark Superligaen" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Denmark Superligaen")

# Row 316
$ws.Range("B316").Value = 5143633
$ws.Range("F316").Value = "AGF Aarhus"
$ws.Range("G316").Value = "Viborg"
$ws.Range("H316").Value = 3
$ws.Range("I316").Value = 1
$ws.Range("J316").Value = "H"
$ws.Range("K316").Value = 2.5
$ws.Range("L316").Value = 3.4
$ws.Range("M316").Value = 2.7
$ws.Range("N316").Value = 2.875
$ws.Range("O316").Value = 3.4
$ws.Range("P316").Value = 2.5
$ws.Range("Q316").Value = 0.25
$ws.Range("R316").Value = 1.79
$ws.Range("S316").Value = 2.11
$ws.Range("T316").Value = 2.5
$ws.Range("U316").Value = 2.025
$ws.Range("V316").Value = 1.825
$ws.Range("W316").Value = 1.875
$ws.Range("X316").Value = -1
$ws.Range("Y316").Value = -1
$ws.Range("Z316").Value = 0.79
$ws.Range("AA316").Value = -1
$ws.Range("AB316").Value = 1.025
$ws.Range("AC316").Value = -1

# Row 317
$ws.Range("B317").Value = 5143634
$ws.Range("F317").Value = "Randers FC"
$ws.Range("G317").Value = "Odense BK"
$ws.Range("H317").Value = 2
$ws.Range("I317").Value = 2
$ws.Range("J317").Value = "D"
$ws.Range("K317").Value = 2
$ws.Range("L317").Value = 3.5
$ws.Range("M317").Value = 3.4
$ws.Range("N317").Value = 2.15
$ws.Range("O317").Value = 3.8
$ws.Range("P317").Value = 3.25
$ws.Range("Q317").Value = -0.25
$ws.Range("R317").Value = 1.88
$ws.Range("S317").Value = 2.02
$ws.Range("T317").Value = 2.75
$ws.Range("U317").Value = 1.875
$ws.Range("V317").Value = 1.975
$ws.Range("W317").Value = -1
$ws.Range("X317").Value = 2.8
$ws.Range("Y317").Value = -1
$ws.Range("Z317").Value = -0.5
$ws.Range("AA317").Value = 0.51
$ws.Range("AB317").Value = 0.875
$ws.Range("AC317").Value = -1

# Row 418
$ws.Range("B418").Value = 5143716
$ws.Range("F418").Value = "Randers FC"
$ws.Range("G418").Value = "Lyngby"
$ws.Range("H418").Value = 1
$ws.Range("I418").Value = 0
$ws.Range("J418").Value = "H"
$ws.Range("K418").Value = 1.85
$ws.Range("L418").Value = 3.75
$ws.Range("M418").Value = 3.75
$ws.Range("N418").Value = 1.909
$ws.Range("O418").Value = 3.75
$ws.Range("P418").Value = 3.8
$ws.Range("Q418").Value = -0.5
$ws.Range("R418").Value = 1.875
$ws.Range("S418").Value = 1.975
$ws.Range("T418").Value = 2.5
$ws.Range("U418").Value = 1.875
$ws.Range("V418").Value = 1.975
$ws.Range("W418").Value = 0.909
$ws.Range("X418").Value = -1
$ws.Range("Y418").Value = -1
$ws.Range("Z418").Value = 0.875
$ws.Range("AA418").Value = -1
$ws.Range("AB418").Value = -1
$ws.Range("AC418").Value = 0.9750000000000001

# Row 419
$ws.Range("B419").Value = 5143701
$ws.Range("F419").Value = "AC Horsens"
$ws.Range("G419").Value = "Viborg"
$ws.Range("H419").Value = 0
$ws.Range("I419").Value = 3
$ws.Range("J419").Value = "A"
$ws.Range("K419").Value = 3.5
$ws.Range("L419").Value = 3.6
$ws.Range("M419").Value = 1.909
$ws.Range("N419").Value = 3.4
$ws.Range("O419").Value = 3.5
$ws.Range("P419").Value = 2.1
$ws.Range("Q419").Value = 0.25
$ws.Range("R419").Value = 2.025
$ws.Range("S419").Value = 1.825
$ws.Range("T419").Value = 2.25
$ws.Range("U419").Value = 1.85
$ws.Range("V419").Value = 2
$ws.Range("W419").Value = -1
$ws.Range("X419").Value = -1
$ws.Range("Y419").Value = 1.1
$ws.Range("Z419").Value = -1
$ws.Range("AA419").Value = 0.825
$ws.Range("AB419").Value = 0.8500000000000001
$ws.Range("AC419").Value = -1

# Row 438
$ws.Range("B438").Value = 5143710
$ws.Range("F438").Value = "AaB"
$ws.Range("G438").Value = "Randers FC"
$ws.Range("H438").Value = 0
$ws.Range("I438").Value = 1
$ws.Range("J438").Value = "A"
$ws.Range("K438").Value = 2.55
$ws.Range("L438").Value = 3.4
$ws.Range("M438").Value = 2.6
$ws.Range("N438").Value = 2.375
$ws.Range("O438").Value = 3.5
$ws.Range("P438").Value = 2.8
$ws.Range("Q438").Value = 0
$ws.Range("R438").Value = 1.8
$ws.Range("S438").Value = 2.05
$ws.Range("T438").Value = 2.5
$ws.Range("U438").Value = 2.025
$ws.Range("V438").Value = 1.825
$ws.Range("W438").Value = -1
$ws.Range("X438").Value = -1
$ws.Range("Y438").Value = 1.8
$ws.Range("Z438").Value = -1
$ws.Range("AA438").Value = 1.05
$ws.Range("AB438").Value = -1
$ws.Range("AC438").Value = 0.825

# Row 439
$ws.Range("B439").Value = 5143712
$ws.Range("F439").Value = "Lyngby"
$ws.Range("G439").Value = "AC Horsens"
$ws.Range("H439").Value = 1
$ws.Range("I439").Value = 1
$ws.Range("J439").Value = "D"
$ws.Range("K439").Value = 2.2
$ws.Range("L439").Value = 3.4
$ws.Range("M439").Value = 3.2
$ws.Range("N439").Value = 2.1
$ws.Range("O439").Value = 3.4
$ws.Range("P439").Value = 3.75
$ws.Range("Q439").Value = -0.25
$ws.Range("R439").Value = 1.8
$ws.Range("S439").Value = 2.05
$ws.Range("T439").Value = 2.25
$ws.Range("U439").Value = 1.85
$ws.Range("V439").Value = 2
$ws.Range("W439").Value = -1
$ws.Range("X439").Value = 2.4
$ws.Range("Y439").Value = -1
$ws.Range("Z439").Value = -0.5
$ws.Range("AA439").Value = 0.5249999999999999
$ws.Range("AB439").Value = -0.5
$ws.Range("AC439").Value = 0.5

# Row 440
$ws.Range("B440").Value = 5143709
$ws.Range("F440").Value = "AGF Aarhus"
$ws.Range("G440").Value = "Odense BK"
$ws.Range("H440").Value = 1
$ws.Range("I440").Value = 0
$ws.Range("J440").Value = "H"
$ws.Range("K440").Value = 1.95
$ws.Range("L440").Value = 3.5
$ws.Range("M440").Value = 3.8
$ws.Range("N440").Value = 1.833
$ws.Range("O440").Value = 3.75
$ws.Range("P440").Value = 4
$ws.Range("Q440").Value = -0.75
$ws.Range("R440").Value = 2.05
$ws.Range("S440").Value = 1.8
$ws.Range("T440").Value = 2.5
$ws.Range("U440").Value = 1.9
$ws.Range("V440").Value = 1.95
$ws.Range("W440").Value = 0.833
$ws.Range("X440").Value = -1
$ws.Range("Y440").Value = -1
$ws.Range("Z440").Value = 0.5249999999999999
$ws.Range("AA440").Value = -0.5
$ws.Range("AB440").Value = -1
$ws.Range("AC440").Value = 0.95

# Row 460
$ws.Range("B460").Value = 6437824
$ws.Range("F460").Value = "Viborg"
$ws.Range("G460").Value = "FC Nordsjaelland"
$ws.Range("H460").Value = 1
$ws.Range("I460").Value = 0
$ws.Range("J460").Value = "H"
$ws.Range("K460").Value = 2.7
$ws.Range("L460").Value = 3.6
$ws.Range("M460").Value = 2.4
$ws.Range("N460").Value = 2.6
$ws.Range("O460").Value = 3.5
$ws.Range("P460").Value = 2.55
$ws.Range("Q460").Value = 0
$ws.Range("R460").Value = 1.95
$ws.Range("S460").Value = 1.9
$ws.Range("T460").Value = 2.5
$ws.Range("U460").Value = 1.9
$ws.Range("V460").Value = 1.95
$ws.Range("W460").Value = 1.6
$ws.Range("X460").Value = -1
$ws.Range("Y460").Value = -1
$ws.Range("Z460").Value = 0.95
$ws.Range("AA460").Value = -1
$ws.Range("AB460").Value = -1
$ws.Range("AC460").Value = 0.95

# Row 461
$ws.Range("B461").Value = 6471200
$ws.Range("F461").Value = "AC Horsens"
$ws.Range("G461").Value = "Midtjylland"
$ws.Range("H461").Value = 0
$ws.Range("I461").Value = 2
$ws.Range("J461").Value = "A"
$ws.Range("K461").Value = 4.5
$ws.Range("L461").Value = 3.75
$ws.Range("M461").Value = 1.714
$ws.Range("N461").Value = 5
$ws.Range("O461").Value = 4
$ws.Range("P461").Value = 1.666
$ws.Range("Q461").Value = 0.75
$ws.Range("R461").Value = 1.975
$ws.Range("S461").Value = 1.875
$ws.Range("T461").Value = 2.5
$ws.Range("U461").Value = 1.825
$ws.Range("V461").Value = 2.025
$ws.Range("W461").Value = -1
$ws.Range("X461").Value = -1
$ws.Range("Y461").Value = 0.6659999999999999
$ws.Range("Z461").Value = -1
$ws.Range("AA461").Value = 0.875
$ws.Range("AB461").Value = -1
$ws.Range("AC461").Value = 1.025

# Row 466
$ws.Range("B466").Value = 6478387
$ws.Range("F466").Value = "AaB"
$ws.Range("G466").Value = "Lyngby"
$ws.Range("H466").Value = 1
$ws.Range("I466").Value = 0
$ws.Range("J466").Value = "H"
$ws.Range("K466").Value = 1.9
$ws.Range("L466").Value = 3.8
$ws.Range("M466").Value = 3.5
$ws.Range("N466").Value = 1.727
$ws.Range("O466").Value = 4
$ws.Range("P466").Value = 4.5
$ws.Range("Q466").Value = -0.75
$ws.Range("R466").Value = 1.98
$ws.Range("S466").Value = 1.92
$ws.Range("T466").Value = 2.75
$ws.Range("U466").Value = 2
$ws.Range("V466").Value = 1.85
$ws.Range("W466").Value = 0.7270000000000001
$ws.Range("X466").Value = -1
$ws.Range("Y466").Value = -1
$ws.Range("Z466").Value = 0.49
$ws.Range("AA466").Value = -0.5
$ws.Range("AB466").Value = -1
$ws.Range("AC466").Value = 0.8500000000000001

# Row 467
$ws.Range("B467").Value = 6437825
$ws.Range("F467").Value = "Viborg"
$ws.Range("G467").Value = "Randers FC"
$ws.Range("H467").Value = 3
$ws.Range("I467").Value = 1
$ws.Range("J467").Value = "H"
$ws.Range("K467").Value = 1.9
$ws.Range("L467").Value = 3.6
$ws.Range("M467").Value = 3.6
$ws.Range("N467").Value = 1.909
$ws.Range("O467").Value = 3.6
$ws.Range("P467").Value = 4
$ws.Range("Q467").Value = -0.5
$ws.Range("R467").Value = 1.875
$ws.Range("S467").Value = 1.975
$ws.Range("T467").Value = 2.5
$ws.Range("U467").Value = 2
$ws.Range("V467").Value = 1.85
$ws.Range("W467").Value = 0.909
$ws.Range("X467").Value = -1
$ws.Range("Y467").Value = -1
$ws.Range("Z467").Value = 0.875
$ws.Range("AA467").Value = -1
$ws.Range("AB467").Value = 1
$ws.Range("AC467").Value = -1

# Row 557
$ws.Range("B557").Value = 6779640
$ws.Range("F557").Value = "Vejle"
$ws.Range("G557").Value = "FC Nordsjaelland"
$ws.Range("H557").Value = 0
$ws.Range("I557").Value = 0
$ws.Range("J557").Value = "D"
$ws.Range("K557").Value = 4.333
$ws.Range("L557").Value = 3.8
$ws.Range("M557").Value = 1.727
$ws.Range("N557").Value = 5
$ws.Range("O557").Value = 4
$ws.Range("P557").Value = 1.666
$ws.Range("Q557").Value = 0.75
$ws.Range("R557").Value = 2
$ws.Range("S557").Value = 1.85
$ws.Range("T557").Value = 2.5
$ws.Range("U557").Value = 1.85
$ws.Range("V557").Value = 2
$ws.Range("W557").Value = -1
$ws.Range("X557").Value = 3
$ws.Range("Y557").Value = -1
$ws.Range("Z557").Value = 1
$ws.Range("AA557").Value = -1
$ws.Range("AB557").Value = -1
$ws.Range("AC557").Value = 1

# Row 558
$ws.Range("B558").Value = 6779638
$ws.Range("F558").Value = "Randers FC"
$ws.Range("G558").Value = "Silkeborg IF"
$ws.Range("H558").Value = 1
$ws.Range("I558").Value = 0
$ws.Range("J558").Value = "H"
$ws.Range("K558").Value = 3
$ws.Range("L558").Value = 3.6
$ws.Range("M558").Value = 2.15
$ws.Range("N558").Value = 3.2
$ws.Range("O558").Value = 3.6
$ws.Range("P558").Value = 2.15
$ws.Range("Q558").Value = 0.25
$ws.Range("R558").Value = 1.925
$ws.Range("S558").Value = 1.925
$ws.Range("T558").Value = 2.5
$ws.Range("U558").Value = 1.95
$ws.Range("V558").Value = 1.9
$ws.Range("W558").Value = 2.2
$ws.Range("X558").Value = -1
$ws.Range("Y558").Value = -1
$ws.Range("Z558").Value = 0.925
$ws.Range("AA558").Value = -1
$ws.Range("AB558").Value = -1
$ws.Range("AC558").Value = 0.8999999999999999

# Row 563
$ws.Range("B563").Value = 6779645
$ws.Range("F563").Value = "Vejle"
$ws.Range("G563").Value = "Hvidovre IF"
$ws.Range("H563").Value = 3
$ws.Range("I563").Value = 1
$ws.Range("J563").Value = "H"
$ws.Range("K563").Value = 1.833
$ws.Range("L563").Value = 3.6
$ws.Range("M563").Value = 4.2
$ws.Range("N563").Value = 1.8
$ws.Range("O563").Value = 3.6
$ws.Range("P563").Value = 4.5
$ws.Range("Q563").Value = -0.5
$ws.Range("R563").Value = 1.825
$ws.Range("S563").Value = 2.025
$ws.Range("T563").Value = 2.5
$ws.Range("U563").Value = 2.025
$ws.Range("V563").Value = 1.825
$ws.Range("W563").Value = 0.8
$ws.Range("X563").Value = -1
$ws.Range("Y563").Value = -1
$ws.Range("Z563").Value = 0.825
$ws.Range("AA563").Value = -1
$ws.Range("AB563").Value = 1.025
$ws.Range("AC563").Value = -1

# Row 564
$ws.Range("B564").Value = 6779644
$ws.Range("F564").Value = "FC Nordsjaelland"
$ws.Range("G564").Value = "Odense BK"
$ws.Range("H564").Value = 0
$ws.Range("I564").Value = 1
$ws.Range("J564").Value = "A"
$ws.Range("K564").Value = 1.5
$ws.Range("L564").Value = 4.2
$ws.Range("M564").Value = 6
$ws.Range("N564").Value = 1.333
$ws.Range("O564").Value = 5
$ws.Range("P564").Value = 9
$ws.Range("Q564").Value = -1.5
$ws.Range("R564").Value = 1.85
$ws.Range("S564").Value = 2
$ws.Range("T564").Value = 3.5
$ws.Range("U564").Value = 2.025
$ws.Range("V564").Value = 1.825
$ws.Range("W564").Value = -1
$ws.Range("X564").Value = -1
$ws.Range("Y564").Value = 8
$ws.Range("Z564").Value = -1
$ws.Range("AA564").Value = 1
$ws.Range("AB564").Value = -1
$ws.Range("AC564").Value = 0.825
